# Insert a new daily-price record above row 214 (Vega Central Mapocho de
# Santiago - Poroto granado). This pushes the existing rows 214-308 down to
# 215-309 and grows the used range to A1:R309.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("214:214").Insert()

$ws.Range("A214").Value = 9
$ws.Range("B214").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C214").Value = "Metropolitana"
$ws.Range("D214").Value = 44875
$ws.Range("E214").Value = 13
$ws.Range("F214").Value = 100112030
$ws.Range("G214").Value = "Poroto granado"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 35
$ws.Range("K214").Value = 48000
$ws.Range("L214").Value = 48000
$ws.Range("M214").Value = 48000
$ws.Range("N214").Value = "$/malla 25 kilos"
$ws.Range("O214").Value = "Perú"
$ws.Range("P214").Value = 1920
$ws.Range("Q214").Value = 25
$ws.Range("R214").Value = "Hortaliza"
